# Updated symbol list refresh (crypto price/volume/hour snapshot).
# Every data row's "Hora" (hour) column moves from 13 -> 14, and most rows'
# Price (D) / Volume(1h) (E) columns pick up a refreshed quote. Rows whose
# coin had no quote at all ("--" / "--%") only get the hour bump.
#
# The sheet stores these as plain text (e.g. "330.26", "-0.31%", "14"), not
# numbers, so each value is written with a leading apostrophe to force
# Excel to keep it as text instead of auto-converting it to a number /
# percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2;  D='330.26';       E='-0.31%' }
    @{ Row=3;  D='43.98';        E='5.54%' }
    @{ Row=4;  D='5.600';        E='-1.60%' }
    @{ Row=5;  D='0.08198';      E='-1.85%' }
    @{ Row=6;  D='8.776';        E='-0.52%' }
    @{ Row=7;  D='4.422';        E='-2.55%' }
    @{ Row=8;  D='1.918';        E='-5.41%' }
    @{ Row=9;  D='2.859';        E='-4.34%' }
    @{ Row=10; D='0.9430';       E='1.40%' }
    @{ Row=11; D='0.1211';       E='-6.53%' }
    @{ Row=12; D='0.1935';       E='-1.83%' }
    @{ Row=13; D='0.09851';      E='4.55%' }
    @{ Row=14; D='0.04426';      E='12.82%' }
    @{ Row=15; D=$null;          E='0.75%' }
    @{ Row=16; D='0.001279';     E='-2.12%' }
    @{ Row=17; D='0.005982';     E='-3.32%' }
    @{ Row=18; D='3.501';        E='1.66%' }
    @{ Row=19; D=$null;          E=$null }
    @{ Row=20; D='8.725';        E='5.05%' }
    @{ Row=21; D=$null;          E='0.61%' }
    @{ Row=22; D='0.2522';       E='1.55%' }
    @{ Row=23; D='0.04397';      E='-0.25%' }
    @{ Row=24; D='0.001239';     E='-0.67%' }
    @{ Row=25; D='0.004315';     E='-1.83%' }
    @{ Row=26; D='0.0001234';    E='2.81%' }
    @{ Row=27; D='0.0004002';    E='31.42%' }
    @{ Row=28; D=$null;          E=$null }
    @{ Row=29; D=$null;          E=$null }
    @{ Row=30; D=$null;          E=$null }
    @{ Row=31; D=$null;          E=$null }
    @{ Row=32; D=$null;          E=$null }
    @{ Row=33; D=$null;          E=$null }
    @{ Row=34; D=$null;          E=$null }
    @{ Row=35; D=$null;          E=$null }
    @{ Row=36; D=$null;          E=$null }
    @{ Row=37; D=$null;          E=$null }
    @{ Row=38; D=$null;          E=$null }
    @{ Row=39; D='0.02847';      E='1.69%' }
    @{ Row=40; D='0.05725';      E='3.30%' }
    @{ Row=41; D='0.007919';     E=$null }
    @{ Row=42; D='0.009800';     E='9.71%' }
    @{ Row=43; D='0.1416';       E='-1.51%' }
    @{ Row=44; D='0.002099';     E='-1.94%' }
    @{ Row=45; D='0.009787';     E='-16.80%' }
    @{ Row=46; D='0.00007309';   E='4.00%' }
    @{ Row=47; D='0.00000000752';E='0.31%' }
    @{ Row=48; D='0.003287';     E='-6.13%' }
    @{ Row=49; D='0.002277';     E='-0.02%' }
    @{ Row=50; D='0.00002107';   E='0.31%' }
    @{ Row=51; D='0.0002007';    E='0.31%' }
)

foreach ($r in $rows) {
    if ($null -ne $r.D) {
        $ws.Cells.Item($r.Row, 4).Value = "'" + $r.D
    }
    if ($null -ne $r.E) {
        $ws.Cells.Item($r.Row, 5).Value = "'" + $r.E
    }
    # Every row's "Hora" column advances from 13 to 14.
    $ws.Cells.Item($r.Row, 7).Value = "'14"
}
